$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3077.6538
$ws.Range("I76").Value = 3096.3635
$ws.Range("J76").Value = 3063.9333
$ws.Range("K76").Value = 3096.3635
$ws.Range("L76").Value = 3063.9333
$ws.Range("M76").Value = -2781.3635
$ws.Range("N76").Value = -3693.9333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3077.6538
$ws.Range("I79").Value = 3096.3635
$ws.Range("J79").Value = 3063.9333
$ws.Range("K79").Value = 3096.3635
$ws.Range("L79").Value = 3063.9333
$ws.Range("M79").Value = -2004.3635
$ws.Range("N79").Value = -5247.933300000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4187.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4187.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4187.5
$ws.Range("N86").Value = -6433.5
$ws.Range("M86").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4187.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 4187.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 20937.5
$ws.Range("N89").Value = -32169.5
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 716287.1
$ws.Range("I106").Value = 716287.1
$ws.Range("K106").Value = 716287.1
$ws.Range("M106").Value = -715656.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1975.7693
$ws.Range("I113").Value = 1726.6666
$ws.Range("J113").Value = 2536.25
$ws.Range("K113").Value = 1726.6666
$ws.Range("L113").Value = 2536.25
$ws.Range("M113").Value = 1527.3334
$ws.Range("N113").Value = -9044.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1482.5
$ws.Range("I2").Value = 1149.4
$ws.Range("J2").Value = 2037.6666
$ws.Range("K2").Value = 1149.4
$ws.Range("L2").Value = 2037.6666
$ws.Range("M2").Value = -1036.4
$ws.Range("N2").Value = -2263.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 23811854
$ws.Range("I61").Value = 25002372
$ws.Range("K61").Value = 25002372
$ws.Range("M61").Value = -25002160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 166667800
$ws.Range("I74").Value = 250000500
$ws.Range("J74").Value = 2400
$ws.Range("K74").Value = 250000500
$ws.Range("L74").Value = 2400
$ws.Range("M74").Value = -249999626
$ws.Range("N74").Value = -4148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 166667800
$ws.Range("I77").Value = 250000500
$ws.Range("J77").Value = 2400
$ws.Range("K77").Value = 1250002500
$ws.Range("L77").Value = 12000
$ws.Range("M77").Value = -1249998132
$ws.Range("N77").Value = -20736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1482.5
$ws.Range("I116").Value = 1149.4
$ws.Range("J116").Value = 2037.6666
$ws.Range("K116").Value = 1149.4
$ws.Range("L116").Value = 2037.6666
$ws.Range("M116").Value = 1144.6
$ws.Range("N116").Value = -6625.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 7655.923
$ws.Range("I122").Value = 8236.174000000001
$ws.Range("K122").Value = 24708.522
$ws.Range("M122").Value = -22258.522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 23811854
$ws.Range("I136").Value = 25002372
$ws.Range("K136").Value = 75007116
$ws.Range("M136").Value = -75004566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1482.5
$ws.Range("I3").Value = 1149.4
$ws.Range("J3").Value = 2037.6666
$ws.Range("K3").Value = 1149.4
$ws.Range("L3").Value = 2037.6666
$ws.Range("M3").Value = -1035.4
$ws.Range("N3").Value = -2265.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 20835834
$ws.Range("I86").Value = 2272.4443
$ws.Range("J86").Value = 83336520
$ws.Range("K86").Value = 2272.4443
$ws.Range("L86").Value = 83336520
$ws.Range("M86").Value = -1149.4443
$ws.Range("N86").Value = -83338766

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 20835834
$ws.Range("I89").Value = 2272.4443
$ws.Range("J89").Value = 83336520
$ws.Range("K89").Value = 11362.2215
$ws.Range("L89").Value = 416682600
$ws.Range("M89").Value = -5746.2215
$ws.Range("N89").Value = -416693832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2871.2546
$ws.Range("I105").Value = 1484.4839
$ws.Range("J105").Value = 4662.5
$ws.Range("K105").Value = 1484.4839
$ws.Range("L105").Value = 4662.5
$ws.Range("M105").Value = 262.5161000000001
$ws.Range("N105").Value = -8156.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3075
$ws.Range("I134").Value = 2468.95
$ws.Range("J134").Value = 4176.909
$ws.Range("K134").Value = 7406.849999999999
$ws.Range("L134").Value = 12530.727
$ws.Range("M134").Value = -4871.849999999999
$ws.Range("N134").Value = -17600.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 93860
$ws.Range("J37").Value = 93860
$ws.Range("L37").Value = 281580
$ws.Range("N37").Value = -281804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 598.8823
$ws.Range("I98").Value = 719.5
$ws.Range("J98").Value = 561.7692
$ws.Range("K98").Value = 2158.5
$ws.Range("L98").Value = 1685.3076
$ws.Range("M98").Value = -660.5
$ws.Range("N98").Value = -4681.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15154344
$ws.Range("I80").Value = 37039516
$ws.Range("J80").Value = 3072.6155
$ws.Range("K80").Value = 37039516
$ws.Range("L80").Value = 3072.6155
$ws.Range("M80").Value = -37038518
$ws.Range("N80").Value = -5068.6155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 15154344
$ws.Range("I83").Value = 37039516
$ws.Range("J83").Value = 3072.6155
$ws.Range("K83").Value = 185197580
$ws.Range("L83").Value = 15363.0775
$ws.Range("M83").Value = -185192588
$ws.Range("N83").Value = -25347.0775

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2324.848
$ws.Range("I132").Value = 2058.6667
$ws.Range("J132").Value = 3283.1
$ws.Range("K132").Value = 6176.000100000001
$ws.Range("L132").Value = 9849.299999999999
$ws.Range("M132").Value = -3646.000100000001
$ws.Range("N132").Value = -14909.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 239007.25
$ws.Range("J141").Value = 239007.25
$ws.Range("L141").Value = 239007.25
$ws.Range("N141").Value = -249367.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1524
$ws.Range("I136").Value = 1398.8125
$ws.Range("J136").Value = 1774.375
$ws.Range("K136").Value = 4196.4375
$ws.Range("L136").Value = 5323.125
$ws.Range("M136").Value = -1646.4375
$ws.Range("N136").Value = -10423.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 46683.855
$ws.Range("J138").Value = 46683.855
$ws.Range("L138").Value = 46683.855
$ws.Range("N138").Value = -56963.855
